$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1041632.25
$ws.Range("I33").Value = 1423549.5
$ws.Range("K33").Value = 1423549.5
$ws.Range("M33").Value = -1423320.5
$ws.Range("H100").Value = 54070.15
$ws.Range("I100").Value = 80495
$ws.Range("J100").Value = 4995.4287
$ws.Range("K100").Value = 80495
$ws.Range("L100").Value = 4995.4287
$ws.Range("M100").Value = -79954
$ws.Range("N100").Value = -6077.4287
$ws.Range("H125").Value = 2961.25
$ws.Range("I125").Value = 2970.9092
$ws.Range("K125").Value = 26738.1828
$ws.Range("M125").Value = -24278.1828
$ws.Range("H135").Value = 1205.3226
$ws.Range("I135").Value = 1205.3226
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 10847.9034
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -8312.903399999999
$ws.Range("N135").ClearContents()
$ws.Range("H141").Value = 1942.25
$ws.Range("I141").Value = 2087.3572
$ws.Range("K141").Value = 6262.071599999999
$ws.Range("M141").Value = -1082.071599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3604.153
$ws.Range("I32").Value = 1925.4396
$ws.Range("J32").Value = 25427.428
$ws.Range("K32").Value = 1925.4396
$ws.Range("L32").Value = 25427.428
$ws.Range("M32").Value = -1638.4396
$ws.Range("N32").Value = -26001.428
$ws.Range("H45").Value = 10538
$ws.Range("I45").Value = 11707.333
$ws.Range("K45").Value = 11707.333
$ws.Range("M45").Value = -11330.333
$ws.Range("H61").Value = 4731.122
$ws.Range("I61").Value = 2249.5806
$ws.Range("J61").Value = 12423.9
$ws.Range("K61").Value = 2249.5806
$ws.Range("L61").Value = 12423.9
$ws.Range("M61").Value = -2037.5806
$ws.Range("N61").Value = -12847.9
$ws.Range("H74").Value = 2870.8484
$ws.Range("I74").Value = 2724.7666
$ws.Range("K74").Value = 2724.7666
$ws.Range("M74").Value = -1850.7666
$ws.Range("H77").Value = 2870.8484
$ws.Range("I77").Value = 2724.7666
$ws.Range("K77").Value = 13623.833
$ws.Range("M77").Value = -9255.832999999999
$ws.Range("H93").Value = 35240.332
$ws.Range("J93").Value = 35240.332
$ws.Range("L93").Value = 35240.332
$ws.Range("N93").Value = -40232.332
$ws.Range("H97").Value = 989.5
$ws.Range("I97").Value = 914.1429000000001
$ws.Range("J97").Value = 1868.6666
$ws.Range("K97").Value = 914.1429000000001
$ws.Range("L97").Value = 1868.6666
$ws.Range("M97").Value = -418.1429000000001
$ws.Range("N97").Value = -2860.6666
$ws.Range("H102").Value = 2498.2856
$ws.Range("I102").Value = 2536.6155
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 2536.6155
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -914.6154999999999
$ws.Range("N102").Value = -5244
$ws.Range("H132").Value = 5264.2593
$ws.Range("I132").Value = 4083.1904
$ws.Range("J132").Value = 9398
$ws.Range("K132").Value = 12249.5712
$ws.Range("L132").Value = 28194
$ws.Range("M132").Value = -9719.5712
$ws.Range("N132").Value = -33254
$ws.Range("H136").Value = 4731.122
$ws.Range("I136").Value = 2249.5806
$ws.Range("J136").Value = 12423.9
$ws.Range("K136").Value = 6748.7418
$ws.Range("L136").Value = 37271.7
$ws.Range("M136").Value = -4198.7418
$ws.Range("N136").Value = -42371.7
$ws.Range("H139").Value = 54999.668
$ws.Range("J139").Value = 54999.668
$ws.Range("L139").Value = 54999.668
$ws.Range("N139").Value = -65279.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 529975.1
$ws.Range("I86").Value = 1003340
$ws.Range("J86").Value = 4014.111
$ws.Range("K86").Value = 1003340
$ws.Range("L86").Value = 4014.111
$ws.Range("M86").Value = -1002217
$ws.Range("N86").Value = -6260.111
$ws.Range("H89").Value = 529975.1
$ws.Range("I89").Value = 1003340
$ws.Range("J89").Value = 4014.111
$ws.Range("K89").Value = 5016700
$ws.Range("L89").Value = 20070.555
$ws.Range("M89").Value = -5011084
$ws.Range("N89").Value = -31302.555
$ws.Range("H94").Value = 1374.2941
$ws.Range("I94").Value = 1457.7333
$ws.Range("K94").Value = 1457.7333
$ws.Range("M94").Value = -1006.7333
$ws.Range("H140").Value = 640494.5
$ws.Range("J140").Value = 640494.5
$ws.Range("L140").Value = 640494.5
$ws.Range("N140").Value = -650854.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 6259.6816
$ws.Range("I58").Value = 6892.5
$ws.Range("K58").Value = 6892.5
$ws.Range("M58").Value = -6689.5
$ws.Range("H62").Value = 53745.855
$ws.Range("J62").Value = 8774.333000000001
$ws.Range("L62").Value = 8774.333000000001
$ws.Range("N62").Value = -10022.333
$ws.Range("H65").Value = 53745.855
$ws.Range("J65").Value = 8774.333000000001
$ws.Range("L65").Value = 43871.665
$ws.Range("N65").Value = -50111.665
$ws.Range("H68").Value = 88439.11
$ws.Range("J68").Value = 88439.11
$ws.Range("L68").Value = 88439.11
$ws.Range("N68").Value = -89937.11
$ws.Range("H71").Value = 88439.11
$ws.Range("J71").Value = 88439.11
$ws.Range("L71").Value = 265317.33
$ws.Range("N71").Value = -272805.33
$ws.Range("H82").Value = 44994
$ws.Range("J82").Value = 44994
$ws.Range("L82").Value = 44994
$ws.Range("N82").Value = -45716
$ws.Range("H85").Value = 44994
$ws.Range("J85").Value = 44994
$ws.Range("L85").Value = 44994
$ws.Range("N85").Value = -47490
$ws.Range("H87").Value = 35000
$ws.Range("J87").Value = 35000
$ws.Range("L87").Value = 35000
$ws.Range("N87").Value = -37372
$ws.Range("H90").Value = 35000
$ws.Range("J90").Value = 35000
$ws.Range("L90").Value = 105000
$ws.Range("N90").Value = -116856
$ws.Range("H122").Value = 1883.1818
$ws.Range("I122").Value = 1801.1428
$ws.Range("K122").Value = 5403.428400000001
$ws.Range("M122").Value = -2953.428400000001
$ws.Range("H132").Value = 17739.137
$ws.Range("I132").Value = 11713.312
$ws.Range("K132").Value = 35139.936
$ws.Range("M132").Value = -32609.936
$ws.Range("H133").Value = 47253.445
$ws.Range("I133").Value = 59999.5
$ws.Range("J133").Value = 43611.715
$ws.Range("K133").Value = 59999.5
$ws.Range("L133").Value = 43611.715
$ws.Range("M133").Value = -57469.5
$ws.Range("N133").Value = -48671.715
$ws.Range("H136").Value = 6259.6816
$ws.Range("I136").Value = 6892.5
$ws.Range("K136").Value = 20677.5
$ws.Range("M136").Value = -18127.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1398.1714
$ws.Range("I5").Value = 1091.6666
$ws.Range("J5").Value = 1628.05
$ws.Range("K5").Value = 3274.9998
$ws.Range("L5").Value = 4884.15
$ws.Range("M5").Value = -3162.9998
$ws.Range("N5").Value = -5108.15
$ws.Range("H31").Value = 1200
$ws.Range("I31").Value = 1200
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 3600
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
$ws.Range("M31").Value = -3312
$ws.Range("H122").Value = 25000250
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H131").Value = 1664.0435
$ws.Range("J131").Value = 1782.9678
$ws.Range("L131").Value = 5348.903399999999
$ws.Range("N131").Value = -15428.9034
$ws.Range("H132").Value = 14287422
$ws.Range("I132").Value = 1995.5
$ws.Range("J132").Value = 33334656
$ws.Range("K132").Value = 17959.5
$ws.Range("L132").Value = 300011904
$ws.Range("M132").Value = -15429.5
$ws.Range("N132").Value = -300016964
$ws.Range("H135").Value = 1398.1714
$ws.Range("I135").Value = 1091.6666
$ws.Range("J135").Value = 1628.05
$ws.Range("K135").Value = 9824.999400000001
$ws.Range("L135").Value = 14652.45
$ws.Range("M135").Value = -7289.999400000001
$ws.Range("N135").Value = -19722.45
$ws.Range("H139").Value = 1552.9412
$ws.Range("I139").Value = 1457.2142
$ws.Range("J139").Value = 1999.6666
$ws.Range("K139").Value = 4371.642599999999
$ws.Range("L139").Value = 5998.9998
$ws.Range("M139").Value = 768.3574000000008
$ws.Range("N139").Value = -16278.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 7176.909
$ws.Range("I122").Value = 7176.909
$ws.Range("K122").Value = 21530.727
$ws.Range("M122").Value = -19080.727
$ws.Range("H132").Value = 29211
$ws.Range("I132").Value = 32487.285
$ws.Range("J132").Value = 21566.334
$ws.Range("K132").Value = 97461.855
$ws.Range("L132").Value = 64699.00199999999
$ws.Range("M132").Value = -94931.855
$ws.Range("N132").Value = -69759.00199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5124.75
$ws.Range("I40").Value = 5124.75
$ws.Range("K40").Value = 5124.75
$ws.Range("M40").Value = -4988.75
$ws.Range("H132").Value = 7402.3335
$ws.Range("I132").Value = 7275.04
$ws.Range("J132").Value = 7800.125
$ws.Range("K132").Value = 21825.12
$ws.Range("L132").Value = 23400.375
$ws.Range("M132").Value = -19295.12
$ws.Range("N132").Value = -28460.375
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H136").Value = 2771.585
$ws.Range("I136").Value = 2473.5
$ws.Range("K136").Value = 7420.5
$ws.Range("M136").Value = -4870.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 6385.3335
$ws.Range("I122").Value = 4544.3335
$ws.Range("J122").Value = 8226.333000000001
$ws.Range("K122").Value = 13633.0005
$ws.Range("L122").Value = 24678.999
$ws.Range("M122").Value = -11183.0005
$ws.Range("N122").Value = -29578.999
$ws.Range("H136").Value = 1842.8823
$ws.Range("I136").Value = 971.0909
$ws.Range("K136").Value = 2913.2727
$ws.Range("M136").Value = -363.2727
